$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.6
$ws.Range("G3").Value = 2.18
$ws.Range("L3").Value = 1.35
$ws.Range("P3").Value = 2.24
$ws.Range("Q3").Value = 1.73
$ws.Range("R3").Value = 1.49
$ws.Range("T3").Value = 1.67
$ws.Range("U3").Value = 2.38
$ws.Range("W3").Value = 1.84
$ws.Range("X3").Value = 19.5
$ws.Range("AC3").Value = 8.800000000000001
$ws.Range("AD3").Value = 14.5
$ws.Range("AE3").Value = 38
$ws.Range("AF3").Value = 15.5
$ws.Range("AG3").Value = 11
$ws.Range("AH3").Value = 16
$ws.Range("AI3").Value = 42
$ws.Range("AJ3").Value = 27
$ws.Range("AK3").Value = 22
$ws.Range("AL3").Value = 32
$ws.Range("AM3").Value = 75
$ws.Range("AO3").Value = 30
$ws.Range("F4").Value = 1.68
$ws.Range("H6").Value = 3.85
$ws.Range("I6").Value = 4.5
$ws.Range("F7").Value = 1.91
$ws.Range("H7").Value = 4.2
$ws.Range("P9").Value = 1.33
$ws.Range("Q9").Value = 2.72
$ws.Range("Q10").Value = 3.4
$ws.Range("G11").Value = 1.75
$ws.Range("I11").Value = 8.800000000000001
$ws.Range("Q11").Value = 2
$ws.Range("P12").Value = 2
$ws.Range("T12").Value = 1.83
$ws.Range("U12").Value = 2.16
$ws.Range("Z12").Value = 27
$ws.Range("AJ12").Value = 29
$ws.Range("R13").Value = 1.44
$ws.Range("AF13").Value = 38
$ws.Range("AG13").Value = 19
$ws.Range("H14").Value = 3.7
$ws.Range("N14").Value = 2.88
$ws.Range("O14").Value = 1.51
$ws.Range("P14").Value = 1.62
$ws.Range("T14").Value = 2.1
$ws.Range("U14").Value = 1.87
$ws.Range("AC14").Value = 7.2
$ws.Range("J15").Value = 3.4
$ws.Range("K15").Value = 3.45
$ws.Range("U15").Value = 2.04
$ws.Range("AI15").Value = 65
$ws.Range("J16").Value = 3.35
$ws.Range("P16").Value = 1.7
$ws.Range("T16").Value = 2.08
$ws.Range("U16").Value = 1.91
$ws.Range("AA16").Value = 120
$ws.Range("AN16").Value = 21
$ws.Range("I17").Value = 2.7
